$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Diet"

for ($r = 2; $r -le 10; $r++) {
    $cur = $ws.Cells.Item($r, 4).Text
    if ($cur -eq "Active") {
        $ws.Cells.Item($r, 4).Value = "Feeding"
    } elseif ($cur -eq "Hibernation") {
        $ws.Cells.Item($r, 4).Value = "Fasting"
    }
}

$ws.Range("E1:E1048576").Select()
